$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: "Code chức năng 1(Login)" time 2h -> 4h
$ws.Range("C19").Value = "4h"

# Row 28: "Code các chức năng khác…" time 8h -> 12h
$ws.Range("C28").Value = "12h"
